# Automatic update of files.
# Updates the "Förändrad" (C) date value for all data rows, and re-syncs the
# ordering of the bottom block of rows (7-14) to match the refreshed source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C ("Förändrad") bumped from 46062 to 46063 for every data row ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46063
}

# --- Rows 7-14 got reordered (A, B, G columns) as the underlying list refreshed ---
$rowData = @{
    7  = @{ A = "A 25015-2023"; B = 45085.6989699074;  G = 1.8  }
    8  = @{ A = "A 19922-2025"; B = 45771.63034722222; G = 10.1 }
    9  = @{ A = "A 62884-2021"; B = 44504;              G = 0.8  }
    10 = @{ A = "A 60024-2025"; B = 45992;              G = 1.1  }
    11 = @{ A = "A 3402-2026";  B = 46042.39047453704; G = 5.5  }
    12 = @{ A = "A 14271-2021"; B = 44278;              G = 6.7  }
    13 = @{ A = "A 25634-2025"; B = 45803.59570601852; G = 6    }
    14 = @{ A = "A 28266-2025"; B = 45818.56381944445; G = 1.9  }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 1).Value = $vals.A
    $ws.Cells.Item($r, 2).Value = $vals.B
    $ws.Cells.Item($r, 7).Value = $vals.G
}
